$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Author byline: "Francois Galgani", "Gui-Peng Yang" and "Michel Boufadel"
#    each get split into two runs (as Word does when the spell-checker marks
#    an unrecognised proper noun while you type). The visible text is
#    unchanged, only the run boundary moves. A pure text Find/Replace would
#    collapse the remainder of the paragraph into a single run, so instead
#    we nudge a character-formatting property on (and back off) the target
#    sub-range: that forces the engine to split the run there without
#    touching the text content.
# ---------------------------------------------------------------------------

$r = $d.Content
$r.Find.Execute("Galgani") | Out-Null
$r.Font.Bold = 1
$r.Font.Bold = 0

$r = $d.Content
$r.Find.Execute("Gui") | Out-Null
$r.Font.Bold = 1
$r.Font.Bold = 0

$r = $d.Content
$r.Find.Execute("Michel ") | Out-Null
$r.Font.Bold = 1
$r.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) "ten day period" -> "ten-day period"; drop the stray "average"; and
#    "six month study period" -> "six-month study period".
# ---------------------------------------------------------------------------

$r = $d.Content
$r.Find.Execute("ten day period", $true, $false, $false, $false, $false, `
    $true, 1, $false, "ten-day period", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("exceeding annual external average nutrient loads", $true, $false, $false, $false, $false, `
    $true, 1, $false, "exceeding annual external nutrient loads", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("six month study period", $true, $false, $false, $false, $false, `
    $true, 1, $false, "six-month study period", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Rework the "Our results..." paragraph: the closing sentence ("The
#    results in this paper support the larger conversation...") moves to the
#    front, "close these facilities" -> "close legacy facilities", and
#    "are left to pay the cost of ... " -> "often pay the externalized
#    costs.  Piney Point is only one example of this broader phenomenon."
# ---------------------------------------------------------------------------

$oldPara = "Our results have broad appeal for the environmental concerns of legacy mining activities on coastal resources.  The US state of Florida has historically supported a large fertilizer industry, whereas these mining activities are also a global phenomenon.  Fertilizer production generates a large amount of waste relative to the commercially viable product and many facilities have had insufficient planning to dispose of this waste in an environmentally responsible manner.  Regulatory oversight has also been insufficient to safely and effectively close these facilities.  As a result, environmental resources and taxpayers are left to pay the cost of legacy mining facilities that have not been actively used for many years.  Piney Point is only one example of this broader phenomenon.  The results in this paper support the larger conversation of how insufficient oversight and planning can lead to unintended environmental impacts. "

$newPara = "The results in this paper support the larger conversation of how insufficient oversight and planning can lead to unintended environmental impacts.  The US state of Florida has historically supported a large fertilizer industry, whereas these mining activities are also a global phenomenon.  Fertilizer production generates a large amount of waste relative to the commercially viable product and many facilities have had insufficient planning to dispose of this waste in an environmentally responsible manner.  Regulatory oversight has also been insufficient to safely and effectively close legacy facilities.  As a result, environmental resources and taxpayers often pay the externalized costs.  Piney Point is only one example of this broader phenomenon."

$r = $d.Content
$r.Find.Execute($oldPara, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newPara, 2) | Out-Null

Write-Output "done"
